# Generate Report for Handoff
# Update localization status cells from "In Translation" to "Ready for handoff"
# and refresh the corresponding handoff timestamps on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 20:59:50"

# zh-cn sheet: Status column (C2) and Latest Handoff Datetime (H2).
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 20:59:44"

# de-de sheet: Status column (C2) and Latest Handoff Datetime (H2).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 20:59:50"

# The new status text ("Ready for handoff") is wider than the old one
# ("In Translation"), so Excel re-autofits the status columns on all three
# sheets to the new content's width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
